$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.1034
$ws.Range("E2").Value = 0.03684999999999999
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 198.3
$ws.Range("L2").Value = 0.1909117165687879
$ws.Range("M2").Value = 105.6
$ws.Range("N2").Value = 0.0313157972776608
$ws.Range("O2").Value = 0.5325264750378215
$ws.Range("P2").Value = 73.9
$ws.Range("Q2").Value = 0.02191512707215089
$ws.Range("R2").Value = 0.3726676752395361
$ws.Range("S2").Value = 31.7
$ws.Range("T2").Value = 0.3001893939393939
$ws.Range("U2").Value = 2219.4
$ws.Range("V2").Value = 0.658165534829928
$ws.Range("W2").Value = 0.09784359519613492
$ws.Range("X2").Value = 0.0810769765577076
$ws.Range("Y2").Value = 0.01676661863842732
$ws.Range("Z2").Value = 0.397216007954263
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.06583408901779958
$ws.Range("AC2").Value = -0.06583408901779958
$ws.Range("AD2").Value = 3128.4
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 3128.4
$ws.Range("AG2").Value = 909
$ws.Range("AH2").Value = 0.481255288054765
$ws.Range("AI2").Value = 0.5883215796897038
$ws.Range("AJ2").Value = 0.2123286071336806
$ws.Range("AK2").Value = 0.2934056357122107

# Row 3
$ws.Range("D3").Value = 0.0168
$ws.Range("E3").Value = -0.0182
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 62.2
$ws.Range("L3").Value = 0.2368621477532369
$ws.Range("M3").Value = 35.3
$ws.Range("N3").Value = 0.03607562595809913
$ws.Range("O3").Value = 0.5675241157556269
$ws.Range("P3").Value = 35.3
$ws.Range("Q3").Value = 0.03607562595809913
$ws.Range("R3").Value = 0.5675241157556269
$ws.Range("U3").Value = 694.7
$ws.Range("V3").Value = 0.709964230965764
$ws.Range("W3").Value = 0.07198240944335146
$ws.Range("X3").Value = 0.05746860477704602
$ws.Range("Y3").Value = 0.01451380466630545
$ws.Range("Z3").Value = 1.358158779415567
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.05746860477704602
$ws.Range("AC3").Value = -0.05746860477704602
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = -694.7
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = -2.44785059901339
$ws.Range("AK3").Value = -10.4152923538231

# Row 4
$ws.Range("D4").Value = 0.19
$ws.Range("E4").Value = 0.0919
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 136.1
$ws.Range("L4").Value = 0.1753639994846025
$ws.Range("M4").Value = 70.3
$ws.Range("N4").Value = 0.02936998663101604
$ws.Range("O4").Value = 0.5165319617927994
$ws.Range("P4").Value = 38.6
$ws.Range("Q4").Value = 0.01612633689839572
$ws.Range("R4").Value = 0.2836149889786921
$ws.Range("S4").Value = 31.7
$ws.Range("T4").Value = 0.4509246088193456
$ws.Range("U4").Value = 1524.7
$ws.Range("V4").Value = 0.6369903074866311
$ws.Range("W4").Value = 0.1237047809489184
$ws.Range("X4").Value = 0.1046853483383692
$ws.Range("Y4").Value = 0.0190194326105492
$ws.Range("Z4").Value = 0.3204905847373637
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.07419957325855314
$ws.Range("AC4").Value = -0.07419957325855314
$ws.Range("AD4").Value = 3128.4
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 3128.4
$ws.Range("AG4").Value = 1603.7
$ws.Range("AH4").Value = 0.5665338645418326
$ws.Range("AI4").Value = 0.6866398893790742
$ws.Range("AJ4").Value = 0.4011958071698397
$ws.Range("AK4").Value = 0.5290294913241407

# Remove debt_ebitda (AN) and net_debt_ebitda (AP) columns for these rows
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
